$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d10Text = 'Background
To investigate the acute effects of coronavirus disease (COVID-19) on lifestyle and metabolic parameters in patients with type 1 diabetes (T1D).


Methods
This retrospective cohort study included 34 patients who were admitted to our hospital from April 16 to May 1, 2020. Data on stress levels, sleep duration, exercise, total diet, snacks, and prepared food intake were obtained from the questionnaires.

 Changes in the values of hemoglobin A1c (HbA1c) and body weight from 3 months before the administration of the questionnaire to the time the study questionnaire was administered (pandemic year), and those from 15 months before to 12 months before the administration of the questionnaire (pre-pandemic year) were evaluated.


Results
Increased stress levels and decreased exercise volumes were observed in approximately 60% and 50% of participants during the COVID-19 pandemic, respectively.

 Decreased sleep duration was associated with changes in the body weight for 3 months during pandemic year (r = − 0.40, p = 0.043).

 Furthermore, compared with changes in HbA1c for 3 months during pre-pandemic year, changes in HbA1c during the pandemic year were worse (0.12% [0.33] % during pandemic year vs.

 − 0.09 [0.39] % during pre-pandemic year, p = 0.027).


Conclusions
Many patients experienced stress and exercised less due to the COVID-19 pandemic.

 Glycemic control in patients with T1D was worse than that in the previous year.

 Since the pandemic is currently ongoing, more attention should be paid to stress and lifestyle factor management in patients with T1D.


Supplementary Information
The online version contains supplementary material available at 10.1007/s13340-021-00507-4.
'
$e10Text = '[Yukako%Hosomi%NULL%1,       Chihiro%Munekawa%NULL%1,       Yoshitaka%Hashimoto%y-hashi@koto.kpu-m.ac.jp%1,       Takuro%Okamura%NULL%2,       Takuro%Okamura%NULL%0,       Fuyuko%Takahashi%NULL%1,       Rena%Kawano%NULL%1,       Hanako%Nakajima%NULL%1,       Saori%Majima%NULL%1,       Takafumi%Senmaru%NULL%1,       Naoko%Nakanishi%NULL%1,       Emi%Ushigome%NULL%1,       Masahide%Hamaguchi%NULL%1,       Masahiro%Yamazaki%NULL%1,       Michiaki%Fukui%NULL%1]'

$ws.Range("D10").Value = $d10Text
$ws.Range("E10").Value = $e10Text
